{"js": "// Iteration Conceptual Model Design (3)\n//\n// 1. \"Dienstplaneinstellungen - Ansicht\" body paragraph: rewording about\n//    which data points must be entered (month & year instead of a\n//    \"special Sunday/holiday staffing count\").\n// 2. \"Mitarbeiter anlegen/loeschen - Ansicht\" body paragraph: add a new\n//    \"-Beschaeftigungsbeginn\" field line to the bullet-style field list.\n// 3. \"Abwesenheitsmeldung - Ansicht\" body paragraph: reword to ask for a\n//    time span instead of a single date + shift (drop the shift/dropdown\n//    sentence).\n// 4. \"Tauschanfrage - Ansicht\" body paragraph: reword the first sentence\n//    describing what must be entered for a shift-swap request.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfunction findParagraphByPrefix(items, prefix) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(prefix) === 0) {\n      return items[i];\n    }\n  }\n  throw new Error(\"Paragraph not found for prefix: \" + prefix);\n}\n\n// --- 1. Dienstplaneinstellungen - Ansicht ------------------------------\nconst settingsParagraph = findParagraphByPrefix(\n  paragraphs.items,\n  \"Auf dieser Seite, welche nur f\u00fcr die Stationsleitung zug\u00e4nglich ist,\"\n);\nsettingsParagraph.insertText(\n  \"Auf dieser Seite, welche nur f\u00fcr die Stationsleitung zug\u00e4nglich ist, \" +\n    \"k\u00f6nnen Eckdaten f\u00fcr den vom System zu erstellenden Monatsdienstplan \" +\n    \"festgelegt werden. Zu den einzutragenden Daten geh\u00f6rt zum einen die \" +\n    \"jeweilige Anzahl von Krankenpflegern der vier Schichten (Fr\u00fch, \" +\n    \"Mittel, Sp\u00e4t und Nacht), zum anderen der Monat und das Jahr, f\u00fcr \" +\n    \"welchen ein Dienstplan erstellt werden soll. \u00dcber einen Button \" +\n    \"\\u201cDienstplan generieren\\u201d kann nach dem Eingeben aller \" +\n    \"Eckdaten ein neuer Dienstplan erstellt werden.\",\n  Word.InsertLocation.replace\n);\n\n// --- 2. Mitarbeiter anlegen/l\u00f6schen - Ansicht ---------------------------\nconst rolleResults = body.search(\"-Rolle\", { matchCase: true });\nrolleResults.load(\"items\");\nawait context.sync();\nrolleResults.items[0].insertText(\n  \"\\u000b-Besch\u00e4ftigungsbeginn\",\n  Word.InsertLocation.after\n);\n\n// --- 3. Abwesenheitsmeldung - Ansicht -----------------------------------\nparagraphs.load(\"items/text\");\nawait context.sync();\nconst absenceParagraph = findParagraphByPrefix(\n  paragraphs.items,\n  \"Auf dieser Seite k\u00f6nnen Mitarbeiter Abwesenheitsmeldungen\"\n);\nabsenceParagraph.insertText(\n  \"Auf dieser Seite k\u00f6nnen Mitarbeiter Abwesenheitsmeldungen an das \" +\n    \"System melden. Dazu muss der Zeitraum angegeben werden, in dem der \" +\n    \"Mitarbeiter nicht erscheinen kann. Ein Kommentarfeld bietet die \" +\n    \"M\u00f6glichkeit f\u00fcr Bemerkungen. Durch den Button \\u201cAbwesenheit \" +\n    \"melden\\u201d wird die Abwesenheit in das System gespeist. Also \" +\n    \"Zusatz kann man \u00fcber einen Anhang-Button etwaige Dokumente mit an \" +\n    \"die Abwesenheitsmeldung h\u00e4ngen.\",\n  Word.InsertLocation.replace\n);\n// The paragraph mark itself was bold; the rewritten paragraph no longer is.\nabsenceParagraph.font.bold = false;\n\n// --- 4. Tauschanfrage - Ansicht ------------------------------------------\nconst swapParagraph = findParagraphByPrefix(\n  paragraphs.items,\n  \"Auf dieser Seite k\u00f6nnen Mitarbeiter Tauschanfragen\"\n);\nswapParagraph.insertText(\n  \"Auf dieser Seite k\u00f6nnen Mitarbeiter Tauschanfragen f\u00fcr Schichten \" +\n    \"generieren. Dazu muss das Datum angegeben werden, an dem die zu \" +\n    \"tauschende Schicht liegt. Ein Kommentarfeld bietet die M\u00f6glichkeit \" +\n    \"f\u00fcr Bemerkungen. \u00dcber den Button \\u201cTausch anfragen\\u201d wird \" +\n    \"die Tauschanfrage an das System \u00fcbertragen. Sofern ein Tausch \" +\n    \"vollzogen werden kann, wird der Mitarbeiter diese Schicht nicht \" +\n    \"mehr in seinem Dienstplan haben, und daf\u00fcr eine Schicht an einem \" +\n    \"anderen Tag hinzubekommen.\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Iteration Conceptual Model Design (3)\n#\n# 1. \"Dienstplaneinstellungen - Ansicht\" body paragraph: rewording about\n#    which data points must be entered (month & year instead of a\n#    \"special Sunday/holiday staffing count\").\n# 2. \"Mitarbeiter anlegen/loeschen - Ansicht\" body paragraph: add a new\n#    \"-Beschaeftigungsbeginn\" field line to the bullet-style field list.\n# 3. \"Abwesenheitsmeldung - Ansicht\" body paragraph: reword to ask for a\n#    time span instead of a single date + shift (drop the shift/dropdown\n#    sentence); the paragraph mark is no longer bold.\n# 4. \"Tauschanfrage - Ansicht\" body paragraph: reword the first sentence\n#    describing what must be entered for a shift-swap request.\n\n$d = $word.ActiveDocument\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n# --- 1. Dienstplaneinstellungen - Ansicht -------------------------------\n$d.Content.Find.Execute(\"Auf dieser Seite, welche nur f\u00fcr die Stationsleitung zug\u00e4nglich ist, k\u00f6nnen Eckdaten f\u00fcr den vom System zu erstellenden Monatsdienstplan festgelegt werden. Zu den einzutragenden Daten geh\u00f6rt zum einen die jeweilige Anzahl von Krankenpflegern der vier Schichten (Fr\u00fch, Mittel, Sp\u00e4t und Nacht) und zum anderen die besondere Besetzungsanzahl an Sonn- und Feiertagen. Die letzte Angabe, die get\u00e4tigt werden muss, ist der Monat, f\u00fcr welchen ein Dienstplan erstellt werden soll. \u00dcber einen Button \u201cDienstplan generieren\u201d kann nach dem Eingeben aller Eckdaten ein neuer Dienstplan erstellt werden.\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"Auf dieser Seite, welche nur f\u00fcr die Stationsleitung zug\u00e4nglich ist, k\u00f6nnen Eckdaten f\u00fcr den vom System zu erstellenden Monatsdienstplan festgelegt werden. Zu den einzutragenden Daten geh\u00f6rt zum einen die jeweilige Anzahl von Krankenpflegern der vier Schichten (Fr\u00fch, Mittel, Sp\u00e4t und Nacht), zum anderen der Monat und das Jahr, f\u00fcr welchen ein Dienstplan erstellt werden soll. \u00dcber einen Button \u201cDienstplan generieren\u201d kann nach dem Eingeben aller Eckdaten ein neuer Dienstplan erstellt werden.\", $wdReplaceAll)\n\n# --- 2. Mitarbeiter anlegen/l\u00f6schen - Ansicht ---------------------------\n$d.Content.Find.Execute(\"-Rolle\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"-Rolle^l-Besch\u00e4ftigungsbeginn\", $wdReplaceAll)\n\n# --- 3. Abwesenheitsmeldung - Ansicht ------------------------------------\n$d.Content.Find.Execute(\"Auf dieser Seite k\u00f6nnen Mitarbeiter Abwesenheitsmeldungen an das System melden. Dazu muss das Datum angegeben werden, an welchem der Mitarbeiter nicht erscheinen kann. Danach muss die Schicht angegeben werden, welche durch die Abwesenheit unterbesetzt w\u00e4re. Dies geschieht durch ein Dropdown-Men\u00fc. Ein Kommentarfeld bietet die M\u00f6glichkeit f\u00fcr Bemerkungen. Durch den Button \u201cAbwesenheit melden\u201d wird die Abwesenheit in das System gespeist. Also Zusatz kann man \u00fcber einen Anhang-Button etwaige Dokumente mit an die Abwesenheitsmeldung h\u00e4ngen.\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"Auf dieser Seite k\u00f6nnen Mitarbeiter Abwesenheitsmeldungen an das System melden. Dazu muss der Zeitraum angegeben werden, in dem der Mitarbeiter nicht erscheinen kann. Ein Kommentarfeld bietet die M\u00f6glichkeit f\u00fcr Bemerkungen. Durch den Button \u201cAbwesenheit melden\u201d wird die Abwesenheit in das System gespeist. Also Zusatz kann man \u00fcber einen Anhang-Button etwaige Dokumente mit an die Abwesenheitsmeldung h\u00e4ngen.\", $wdReplaceAll)\n# The paragraph mark itself was bold; the rewritten paragraph no longer is.\n$d.Paragraphs(32).Range.Font.Bold = 0\n\n# --- 4. Tauschanfrage - Ansicht -------------------------------------------\n$d.Content.Find.Execute(\"Auf dieser Seite k\u00f6nnen Mitarbeiter Tauschanfragen f\u00fcr Schichten generieren. Dazu muss zun\u00e4chst das Datum angegeben und \u00fcber ein Dropdown-Men\u00fc die Schicht ausgew\u00e4hlt werden, welche der Mitarbeiter gerne getauscht haben m\u00f6chte. Ein Kommentarfeld bietet die M\u00f6glichkeit f\u00fcr Bemerkungen. \u00dcber den Button \u201cTausch anfragen\u201d wird die Tauschanfrage an das System \u00fcbertragen. Sofern ein Tausch vollzogen werden kann, wird der Mitarbeiter diese Schicht nicht mehr in seinem Dienstplan haben, und daf\u00fcr eine Schicht an einem anderen Tag hinzubekommen.\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"Auf dieser Seite k\u00f6nnen Mitarbeiter Tauschanfragen f\u00fcr Schichten generieren. Dazu muss das Datum angegeben werden, an dem die zu tauschende Schicht liegt. Ein Kommentarfeld bietet die M\u00f6glichkeit f\u00fcr Bemerkungen. \u00dcber den Button \u201cTausch anfragen\u201d wird die Tauschanfrage an das System \u00fcbertragen. Sofern ein Tausch vollzogen werden kann, wird der Mitarbeiter diese Schicht nicht mehr in seinem Dienstplan haben, und daf\u00fcr eine Schicht an einem anderen Tag hinzubekommen.\", $wdReplaceAll)\n"}
